# Update "想去人数" (interested-count) figures in column F across the
# workbook's sheets, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1406
$ws.Range("F4").Value = 25513
$ws.Range("F6").Value = 236
$ws.Range("F7").Value = 567
$ws.Range("F8").Value = 158
$ws.Range("F9").Value = 408
$ws.Range("F11").Value = 338
$ws.Range("F12").Value = 194
$ws.Range("F13").Value = 165
$ws.Range("F14").Value = 40
$ws.Range("F15").Value = 264
$ws.Range("F16").Value = 321
$ws.Range("F17").Value = 42
$ws.Range("F18").Value = 1452
$ws.Range("F19").Value = 144
$ws.Range("F20").Value = 402
$ws.Range("F21").Value = 88

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4496
$ws.Range("F6").Value = 56
$ws.Range("F16").Value = 23

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 4853
$ws.Range("F4").Value = 165

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1406
$ws.Range("F5").Value = 4853
$ws.Range("F6").Value = 165
$ws.Range("F7").Value = 25513
$ws.Range("F9").Value = 4496
$ws.Range("F10").Value = 236
$ws.Range("F13").Value = 567
$ws.Range("F16").Value = 158
$ws.Range("F17").Value = 56
$ws.Range("F18").Value = 56
$ws.Range("F23").Value = 408
$ws.Range("F26").Value = 338
$ws.Range("F27").Value = 194
$ws.Range("F28").Value = 165
$ws.Range("F29").Value = 40
$ws.Range("F31").Value = 264
$ws.Range("F34").Value = 321
$ws.Range("F35").Value = 42
$ws.Range("F37").Value = 1452
$ws.Range("F38").Value = 144
$ws.Range("F39").Value = 23
$ws.Range("F40").Value = 402
$ws.Range("F41").Value = 88
